$wb = $excel.ActiveWorkbook

$wsOptions = $wb.Worksheets.Item("Options")
$wsFutures = $wb.Worksheets.Item("FuturesForwards")
$wsClassify = $wb.Worksheets.Item("ToClassify")
$wsFixedIncome = $wb.Worksheets.Item("Fixed-income")

# --- Options sheet content edits ---
$wsOptions.Range("C3").Value = "La Mc Caulay duration d’une obligation de maturité dans 7ans"
$wsOptions.Range("C7").Value = "Parmis ces graphes suivants, lequel représente pour un long call l'évolution du profit net en fonction du prix du sous-jacent : "
$wsOptions.Range("E7").Value = "images/Question_6_Options/CallOK.JPG"
$wsOptions.Range("F7").Value = "test"

# --- View / selection changes ---
$wsOptions.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsOptions.Range("C9").Select()

$wsFutures.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsFutures.Range("C5").Select()

$wsClassify.Activate()
$wsClassify.Range("E8").Select()

$wsFixedIncome.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
$wsFixedIncome.Range("C3:F3").Select()

$wsOptions.Activate()
